# Template functions & format change
# - Rename the "non_syllabic" feature label to "non syllabic" (disallow the
#   underscore special character in feature descriptions) everywhere it is
#   used on the sheet.
# - Update the saved selection/scroll position (no more functional effect,
#   just reflects where the author left the cursor).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cells = @("B4","B5","B6","B7","B8","B10","B11","B13","B15","B17","B18","B19","B20","B22","B23","B33")
foreach ($addr in $cells) {
    $ws.Range($addr).Value = "non syllabic"
}

$ws.Range("H20").Select() | Out-Null
